$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column B updates (rows 9-17) ---
$ws.Range("B9").Value = 2
$ws.Range("B10").Value = 2
$ws.Range("B11").Value = 2
$ws.Range("B12").Value = 3
$ws.Range("B13").Value = 3
$ws.Range("B14").Value = 3
$ws.Range("B15").Value = 4
$ws.Range("B16").Value = 4
$ws.Range("B17").Value = 4

# --- Column E updates (rows 10-17) ---
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 4
$ws.Range("E13").Value = 6
$ws.Range("E14").Value = 8
$ws.Range("E16").Value = 10
$ws.Range("E17").Value = 0

# --- Column F updates (rows 12, 15) ---
$ws.Range("F12").Value = 5
$ws.Range("F15").Value = 9

# --- J6 / K6 updates ---
$ws.Range("J6").Value = 40
$ws.Range("K6").Formula = "=J6"

# --- Recalculate so dependent formulas (J7:K9, B18, E18, F18) refresh ---
$excel.Calculate()

# --- Column widths for J (10) and K (11) to match best-fit sizing ---
# (input values chosen so the engine's rounded stored width lands as close
# as possible to the authored 11.125 / 9.625 character widths)
$ws.Columns.Item(10).ColumnWidth = 10.4
$ws.Columns.Item(11).ColumnWidth = 8.8

$wb.Save()
